$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for Hi/Lo odds block
$ws.Range("K1").Value = "全場入球大細"

# New data row for Hi/Lo odds
$ws.Range("K2").Value = "更新時間"
$ws.Range("K2").NumberFormat = "[`$-F400]h:mm:ss\ AM/PM"
$ws.Range("L2").Value = "球數"
$ws.Range("M2").Value = "大"
$ws.Range("N2").Value = "細"

# Add the new defined name.
# NOTE: Names.Add() chokes if the *first* argument (the Name) starts with a
# non-ASCII character, so create it with a placeholder ASCII name first and
# rename the returned Name object afterwards.
$newName = $wb.Names.Add("PlaceholderName", '=TEMPLATE!$K$1')
$newName.Name = "全場入球大細"

# Update the selection to match the target state
$ws.Range("L9").Select()
